$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '39.474.88'
$ws.Range("E2").Value = '  +1.94%  '
Set-TextValue $ws.Range("D3") '2.164.01'
$ws.Range("E3").Value = '  +2.97%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue $ws.Range("D5") '227.93'
$ws.Range("E5").Value = '  -0.28%  '
Set-TextValue $ws.Range("D6") '0.623'
$ws.Range("E6").Value = '  +1.12%  '
Set-TextValue $ws.Range("D7") '63.91'
$ws.Range("E7").Value = '  +2.81%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +2.14%  '
Set-TextValue $ws.Range("D10") '0.0855'
$ws.Range("E10").Value = '  +1.67%  '
$ws.Range("E11").Value = '  +0.08%  '
Set-TextValue $ws.Range("D12") '16.16'
$ws.Range("E12").Value = '  +2.03%  '
Set-TextValue $ws.Range("D13") '2.484.90'
$ws.Range("E13").Value = '  +2.96%  '
Set-TextValue $ws.Range("D14") '22.13'
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("E16").Value = '  +0.08%  '
Set-TextValue $ws.Range("D17") '2.175.35'
$ws.Range("E17").Value = '  +3.26%  '
Set-TextValue $ws.Range("D18") '39.488.43'
$ws.Range("E18").Value = '  +1.89%  '
Set-TextValue $ws.Range("D19") '71.96'
$ws.Range("E19").Value = '  +0.19%  '
Set-TextValue $ws.Range("D20") '6.14'
$ws.Range("E20").Value = '  +1.17%  '
$ws.Range("E21").Value = '  +1.31%  '
Set-TextValue $ws.Range("D22") '229.37'
$ws.Range("E22").Value = '  +0.76%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("E25").Value = '  +1.34%  '
Set-TextValue $ws.Range("D26") '172.23'
Set-TextValue $ws.Range("D27") '9.58'
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("E28").Value = '  +1.73%  '
$ws.Range("E29").Value = '  +1.17%  '
Set-TextValue $ws.Range("D30") '19.90'
$ws.Range("E30").Value = '  +2.81%  '
$ws.Range("E31").Value = '  +4.25%  '
$ws.Range("E32").Value = '  +1.24%  '
$ws.Range("E33").Value = '  +1.69%  '
Set-TextValue $ws.Range("D34") '7.11'
$ws.Range("E34").Value = '  +3.72%  '
Set-TextValue $ws.Range("D35") '4.73'
$ws.Range("E35").Value = '  -0.69%  '
$ws.Range("E36").Value = '  -0.25%  '
$ws.Range("E37").Value = '  +0.79%  '
$ws.Range("E38").Value = '  -0.13%  '
$ws.Range("E39").Value = '  -0.17%  '
Set-TextValue $ws.Range("D40") '103.35'
$ws.Range("E40").Value = '  +0.66%  '
$ws.Range("E41").Value = '  +0.84%  '
Set-TextValue $ws.Range("D42") '17.84'
$ws.Range("E42").Value = '  -1.70%  '
Set-TextValue $ws.Range("D43") '1.524.74'
$ws.Range("E43").Value = '  -0.64%  '
$ws.Range("E44").Value = '  +3.39%  '
Set-TextValue $ws.Range("D45") '0.0930'
$ws.Range("E45").Value = '  +2.03%  '
Set-TextValue $ws.Range("D46") '2.82'
$ws.Range("E46").Value = '  +0.73%  '
$ws.Range("E47").Value = '  +5.71%  '
Set-TextValue $ws.Range("D48") '4.28'
$ws.Range("E48").Value = '  +3.50%  '
Set-TextValue $ws.Range("D49") '7.75'
$ws.Range("E49").Value = '  -1.17%  '
Set-TextValue $ws.Range("D50") '2.368.58'
$ws.Range("E50").Value = '  +3.21%  '
$ws.Range("E51").Value = '  -0.38%  '
